$wb = $excel.ActiveWorkbook

# Map of old value -> new value for the F column ("想去人数"), shared between the two
# sheets that list the same expo entries ("展览" starting at row 2, "全部类型" starting at
# row 3, one row offset from the first sheet).
$updates = @{
    1164  = 1165
    1340  = 1342
    296   = 298
    1019  = 1020
    10518 = 10530
    14    = 15
    79    = 80
    281   = 282
    1029  = 1030
    681   = 683
    11975 = 11989
    12396 = 12403
    113   = 114
}

# Sheet "展览": data rows 2-16
$ws1 = $wb.Worksheets.Item("展览")
for ($r = 2; $r -le 16; $r++) {
    $cell = $ws1.Cells.Item($r, 6)  # column F
    $current = $cell.Value2
    if ($null -ne $current -and $updates.ContainsKey([int]$current)) {
        $cell.Value2 = $updates[[int]$current]
    }
}

# Sheet "全部类型": data rows 3-17 (one row offset vs the first sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
for ($r = 3; $r -le 17; $r++) {
    $cell = $ws4.Cells.Item($r, 6)  # column F
    $current = $cell.Value2
    if ($null -ne $current -and $updates.ContainsKey([int]$current)) {
        $cell.Value2 = $updates[[int]$current]
    }
}
